# Update crypto price/volume data cells per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.577.51'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '3.070.82'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.76%  '
$ws.Range('D5').Value = '''592.65'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').Value = '''154.92'
$ws.Range('E6').Value = '  +1.95%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('D9').Value = '3.070.52'
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('E10').Value = '  -1.16%  '
$ws.Range('D11').Value = '''5.92'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('E12').Value = '  -1.99%  '
$ws.Range('D13').Value = '''0.0000237'
$ws.Range('E13').Value = '  -2.20%  '
$ws.Range('D14').Value = '''36.63'
$ws.Range('E14').Value = '  -3.15%  '
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '3.578.90'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').Value = '''7.18'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').Value = '63.472.82'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').Value = '3.075.56'
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('D20').Value = '''479.40'
$ws.Range('E20').Value = '  +2.22%  '
$ws.Range('D21').Value = '''14.45'
$ws.Range('E21').Value = '  -2.79%  '
$ws.Range('D22').Value = '''0.709'
$ws.Range('E22').Value = '  -3.94%  '
$ws.Range('D23').Value = '''7.55'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = '''2.40'
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('D25').Value = '''81.73'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = '''12.85'
$ws.Range('E26').Value = '  -3.17%  '
$ws.Range('D27').Value = '''10.77'
$ws.Range('E27').Value = '  +9.43%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').Value = '''7.64'
$ws.Range('E29').Value = '  +2.92%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''2.21'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('E33').Value = '  -3.87%  '
$ws.Range('D34').Value = '''27.18'
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('D35').Value = '0.0₃0827'
$ws.Range('E35').Value = '  -2.89%  '
$ws.Range('D36').Value = '''1.06'
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('D37').Value = '''6.07'
$ws.Range('E37').Value = '  -1.17%  '
$ws.Range('D38').Value = '''3.28'
$ws.Range('E38').Value = '  -2.82%  '
$ws.Range('D39').Value = '''2.23'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').Value = '''9.23'
$ws.Range('E40').Value = '  -1.36%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '''50.60'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('D42').Value = '''442.00'
$ws.Range('E42').Value = '  -2.35%  '
$ws.Range('D43').Value = '''0.291'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('E44').Value = '  +2.91%  '
$ws.Range('E45').Value = '  -2.31%  '
$ws.Range('D46').Value = '''39.90'
$ws.Range('E46').Value = '  +2.90%  '
$ws.Range('D47').Value = '2.821.54'
$ws.Range('E47').Value = '  -1.00%  '
$ws.Range('D48').Value = '''132.35'
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').Value = '''0.999'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '''25.25'
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('D51').Value = '''2.24'
$ws.Range('E51').Value = '  -1.38%  '
